# HZNPNamedTrade20.xlsx - "traded, fixed issues with the repeater"
#
# A new trade was recorded, so append it as row 10 of the log (the sheet
# previously ran from row 1 header + rows 2-9 of data). The new row
# mirrors the layout/format of the existing trade rows (3-9): dates in
# column A and the IsShortSell flag in column G keep the workbook's
# date-style (s="1"); everything else is a plain number/boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 9's formatting down into row 10 first (this carries the
# date-format style used on A10/G10) so the new cells land in the same
# style slots as the existing rows, then overwrite with the new values.
$ws.Range("A9:I9").Copy() | Out-Null
$ws.Range("A10:I10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(10, 1).Value = 42654.746006944442      # Date
$ws.Cells.Item(10, 2).Value = $true                   # Profitable
$ws.Cells.Item(10, 3).Value = 10022.16                # Principle
$ws.Cells.Item(10, 4).Value = 10009.15                # Start Principle
$ws.Cells.Item(10, 5).Value = 18.870000999999998      # BuyPrice
$ws.Cells.Item(10, 6).Value = 18.920000000000002      # SellPrice
$ws.Cells.Item(10, 7).Value = $false                  # IsShortSell
$ws.Cells.Item(10, 8).Value = 0.26                    # Price Change %
$ws.Cells.Item(10, 9).Value = $false                  # Strong trade

# The sheet's columns are all "best fit" (auto-sized to content) -- with
# the new row in place, re-apply the refreshed best-fit widths the same
# way Excel recalculates them whenever the underlying data changes.
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.333333333333333
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 10.333333333333334
$ws.Columns.Item(5).ColumnWidth = 9
$ws.Columns.Item(6).ColumnWidth = 6.166666666666667
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 13.833333333333334
$ws.Columns.Item(9).ColumnWidth = 11
